$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "ok"
$ws.Range("F2").Value = 0.21399174797286
$ws.Range("H2").Value = 0.825
$ws.Range("J2").Value = 7.75
$ws.Range("K2").Value = 0.606658753564947
$ws.Range("L2").Value = -0.593574370027795
$ws.Range("M2").Value = 1.65718531531839
$ws.Range("N2").Value = 7.82785488470899
$ws.Range("P2").Value = "Unlikely improving"
# Row 3
$ws.Range("D3").Value = $true
$ws.Range("F3").Value = 0.020915297282622
$ws.Range("H3").Value = 0.87719298245614
$ws.Range("K3").Value = -0.074142327685791
$ws.Range("L3").Value = -0.149619290428151
$ws.Range("M3").Value = -0.0138649198239153
$ws.Range("N3").Value = -0.701441132315904
$ws.Range("P3").Value = "Extremely unlikely increasing"
# Row 4
$ws.Range("F4").Value = 0.0140997155703917
$ws.Range("G4").Value = 0.610169491525424
$ws.Range("H4").Value = 0.152542372881356
$ws.Range("P4").Value = "Extremely unlikely improving"
# Row 5
$ws.Range("F5").Value = 0.469683596731774
$ws.Range("G5").Value = 0.867924528301887
$ws.Range("H5").Value = 0.188679245283019
$ws.Range("I5").Value = 3
$ws.Range("P5").Value = "As likely as not improving"
# Row 6
$ws.Range("F6").Value = 0.9999972314077
$ws.Range("G6").Value = 0.423728813559322
$ws.Range("H6").Value = 0.11864406779661
$ws.Range("J6").Value = 0.002
$ws.Range("K6").Value = -0.0005929383116882999
$ws.Range("L6").Value = -0.0008370893812069999
$ws.Range("N6").Value = -29.6469155844156
$ws.Range("P6").Value = "Virtually certain improving"
# Row 7
$ws.Range("F7").Value = 0.830172412503015
$ws.Range("H7").Value = 0.932203389830508
$ws.Range("K7").Value = -0.0206540178571429
$ws.Range("L7").Value = -0.0586258800847643
$ws.Range("M7").Value = 0.0163488968027403
$ws.Range("N7").Value = -5.44960893328307
$ws.Range("P7").Value = "Likely improving"
# Row 8
$ws.Range("F8").Value = 0.0036332507713829
$ws.Range("H8").Value = 0.857142857142857
$ws.Range("J8").Value = 7.505
$ws.Range("K8").Value = -0.0852918956043956
$ws.Range("L8").Value = -0.130446428571428
$ws.Range("M8").Value = -0.0540693737416814
$ws.Range("N8").Value = -1.13646762963885
# Row 9
$ws.Range("F9").Value = 0.8285081831912749
$ws.Range("H9").Value = 0.949152542372881
$ws.Range("J9").Value = 0.389
$ws.Range("K9").Value = -0.0207984515484516
$ws.Range("L9").Value = -0.0558575491646487
$ws.Range("M9").Value = 0.0164771103994961
$ws.Range("N9").Value = -5.34664564227546
$ws.Range("P9").Value = "Likely improving"
# Row 10
$ws.Range("D10").Value = $true
$ws.Range("F10").Value = 0.614457968297824
$ws.Range("H10").Value = 0.6610169491525421
$ws.Range("J10").Value = 0.55
$ws.Range("K10").Value = -0.0041583203085885
$ws.Range("L10").Value = -0.0238892595615511
$ws.Range("M10").Value = 0.0208044514574797
$ws.Range("N10").Value = -0.756058237925181
$ws.Range("P10").Value = "As likely as not improving"
# Row 11
$ws.Range("F11").Value = 0.0809039929491798
$ws.Range("H11").Value = 0.542372881355932
$ws.Range("J11").Value = 0.011
$ws.Range("K11").Value = 0.000909336099585
$ws.Range("M11").Value = 0.0021882191423468
$ws.Range("N11").Value = 8.266691814409651
$ws.Range("P11").Value = "Very unlikely improving"
# Row 12
$ws.Range("F12").Value = 0.0344033944111115
$ws.Range("H12").Value = 0.760869565217391
$ws.Range("J12").Value = 4.525
$ws.Range("K12").Value = 0.30103021978022
$ws.Range("L12").Value = 0.0108190993747977
$ws.Range("M12").Value = 0.614887579117234
$ws.Range("N12").Value = 6.6526015421043
$ws.Range("P12").Value = "Extremely unlikely improving"
# Row 13
$ws.Range("F13").Value = 0.9486359866905399
$ws.Range("H13").Value = 0.815789473684211
$ws.Range("K13").Value = 0.0339860426929392
$ws.Range("L13").Value = -0.0011867593108791
$ws.Range("M13").Value = 0.0716331106895097
$ws.Range("N13").Value = 0.324759127500614
$ws.Range("P13").Value = "Very likely increasing"
# Row 14
$ws.Range("F14").Value = 0.997852443814679
$ws.Range("G14").Value = 0.444444444444444
$ws.Range("H14").Value = 0.111111111111111
$ws.Range("L14").Value = -0.0002305159612545
$ws.Range("P14").Value = "Virtually certain improving"
# Row 15
$ws.Range("F15").Value = 0.655106009894089
$ws.Range("G15").Value = 0.87962962962963
$ws.Range("H15").Value = 0.185185185185185
$ws.Range("P15").Value = "As likely as not improving"
# Row 16
$ws.Range("F16").Value = 0.999977315254334
$ws.Range("G16").Value = 0.341880341880342
$ws.Range("H16").Value = 0.102564102564103
$ws.Range("K16").Value = -0.0001843765774861
$ws.Range("L16").Value = -0.0003087489433643
$ws.Range("N16").Value = -6.14588591620394
$ws.Range("P16").Value = "Virtually certain improving"
# Row 17
$ws.Range("F17").Value = 0.99850221280582
$ws.Range("H17").Value = 0.957264957264957
$ws.Range("J17").Value = 0.471
$ws.Range("K17").Value = -0.028461038961039
$ws.Range("L17").Value = -0.0457008155828836
$ws.Range("M17").Value = -0.0127603017260923
$ws.Range("N17").Value = -6.04268343121846
# Row 18
$ws.Range("F18").Value = 0.262615025872859
$ws.Range("H18").Value = 0.6754385964912279
$ws.Range("J18").Value = 7.51
$ws.Range("K18").Value = -0.0099911121848508
$ws.Range("L18").Value = -0.0225772664835164
$ws.Range("M18").Value = 0.0109638679094912
$ws.Range("N18").Value = -0.133037445870184
$ws.Range("P18").Value = "Unlikely increasing"
# Row 19
$ws.Range("D19").Value = $true
$ws.Range("F19").Value = 0.998456105826092
$ws.Range("H19").Value = 0.965811965811966
$ws.Range("J19").Value = 0.474
$ws.Range("K19").Value = -0.0304569488817891
$ws.Range("L19").Value = -0.0433948162715718
$ws.Range("M19").Value = -0.0143037716420726
$ws.Range("N19").Value = -6.42551664172767
# Row 20
$ws.Range("F20").Value = 0.9978706440191289
$ws.Range("H20").Value = 0.641025641025641
$ws.Range("J20").Value = 0.62
$ws.Range("K20").Value = -0.0202033274233754
$ws.Range("L20").Value = -0.0361936866357975
$ws.Range("M20").Value = -0.0080860875667559
$ws.Range("N20").Value = -3.25860119731861
$ws.Range("P20").Value = "Virtually certain improving"
# Row 21
$ws.Range("F21").Value = 0.0474326378999681
$ws.Range("G21").Value = 0.0170940170940171
$ws.Range("H21").Value = 0.367521367521368
$ws.Range("K21").Value = 0.0003285420891602
$ws.Range("M21").Value = 0.0007182648135007
$ws.Range("N21").Value = 3.28542089160256
# Row 22
$ws.Range("F22").Value = 0.889664319040077
$ws.Range("K22").Value = 0.028274852246634
$ws.Range("L22").Value = -0.0423542572185553
$ws.Range("M22").Value = 0.0557937693963127
$ws.Range("N22").Value = 6.17354852546593
$ws.Range("P22").Value = "Likely improving"
# Row 23
$ws.Range("F23").Value = 0.5
$ws.Range("J23").Value = 109.57
$ws.Range("K23").Value = 0.118744992569618
$ws.Range("M23").Value = 4.82506037385531
$ws.Range("N23").Value = 0.108373635638968
$ws.Range("P23").Value = "As likely as not improving"
# Row 24
$ws.Range("F24").Value = 0.95679463351315
$ws.Range("J24").Value = 6.258
$ws.Range("K24").Value = 0.345918038850152
$ws.Range("L24").Value = 0.0324362280233659
$ws.Range("N24").Value = 5.52761327660837
$ws.Range("P24").Value = "Extremely likely improving"
# Row 25
$ws.Range("F25").Value = 0.105248850027669
$ws.Range("J25").Value = 0.5155
$ws.Range("K25").Value = -0.0148502740977615
$ws.Range("L25").Value = -0.0302033158028902
$ws.Range("M25").Value = 0.0016890463661578
$ws.Range("N25").Value = -2.88075152235917
$ws.Range("P25").Value = "Unlikely improving"
# Row 26
$ws.Range("F26").Value = 0.237137175383665
$ws.Range("J26").Value = 109.285
$ws.Range("K26").Value = -0.75369047619048
$ws.Range("L26").Value = -2.13417463286886
$ws.Range("M26").Value = 0.842346024703969
$ws.Range("N26").Value = -0.689655923677064
$ws.Range("P26").Value = "Unlikely improving"
# Row 27
$ws.Range("F27").Value = 0.570986171506239
$ws.Range("J27").Value = 6.2935
$ws.Range("K27").Value = 0.0241560768698061
$ws.Range("L27").Value = -0.145933780280412
$ws.Range("M27").Value = 0.273257474708751
$ws.Range("N27").Value = 0.383825802332662

Write-Output "applied changes"